# "adding averages and more checks"
#
# 1) Training Dashboard: PERIOD TO EXPIRE (col H) drops by 8 days and
#    LAST UPDATE (col I) moves from 08-Sep-2025 to 16-Sep-2025 for every
#    data row (3-23).
# 2) Exam Dashboard: COMMENTS column (E) gets a wider column and the
#    per-row remark text becomes more descriptive ("date is valid"
#    instead of "OK") for rows 3-4.
# 3) Cosmetic: the dark-blue header banner text turns white, and the big
#    title banner loses its oversized 14pt font (now matches the header
#    banner's bold/white look).

$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: refresh PERIOD TO EXPIRE / LAST UPDATE --------
for ($row = 3; $row -le 23; $row++) {
    $periodCell = $trainingWs.Cells.Item($row, 8)   # column H
    $periodCell.Value = $periodCell.Value2 - 8

    $updateCell = $trainingWs.Cells.Item($row, 9)   # column I
    $updateCell.Value = "'16-Sep-2025"
}

# --- Exam Dashboard: widen COMMENTS column, reword the remarks ---------
$examWs.Columns.Item(5).ColumnWidth = 14.17          # raw width 10 -> 15

$examWs.Range("E3").Value = "date is valid"
$examWs.Range("E4").Value = "date is valid"

# --- Header banner + title font touch-up -------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215
}

$trainingWs.Range("A2:K2").Font.Color = 16777215
$examWs.Range("A2:G2").Font.Color = 16777215
